$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (SNOW): description text loses " Normative" segment.
$ws.Range("C11").Value = "Type=2  Mode=42  LongTime=0.5  ReliabilityFactor=1.4"

# Row 15 (Wind+Y Dynamic): parameter tweak "21 6 10" -> "21 20 10".
$ws.Range("C15").Value = "Type=2  Mode=7  Normative ReliabilityFactor=1.4 21 20 10 1 3 0 0 2 95 6 1 0 0.3 2 0"

# Row 16 (Wind-Y Dynamic): parameter tweak "21 6 10" -> "21 20 11".
$ws.Range("C16").Value = "Type=2  Mode=7  Normative ReliabilityFactor=1.4 21 20 11 1 3 0 0 2 95 6 1 0 0.3 2 0"

# Update the active selection to C19 as recorded in the saved view state.
$ws.Range("C19").Select()
